$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row before row 37 for the new event "杭州·第六届华盟次元动漫嘉年华"
$ws1.Rows.Item(37).Insert()

# Re-apply the index-column style (bold/border/center) to new A37 by copying format from A36
$ws1.Range("A36").Copy()
$ws1.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws1.Range("A37").Value = 36

# Fill new-row content (the new event)
$ws1.Range("B37").NumberFormat = "@"
$ws1.Range("B37").Value = "2024-10-01"
$ws1.Range("C37").Value = "杭州·第六届华盟次元动漫嘉年华"
$ws1.Range("D37").Value = "创意路1号 中国智谷富春园区"
$ws1.Range("E37").Value = "2024.10.01 10:00-10.02 17:00"
$ws1.Range("F37").Value = 4
$ws1.Range("G37").Value = 60
$ws1.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=89966"
$ws1.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202407/d0ryUws41721962610997.jpeg"

# Old row 37 (鸢飞鱼跃) is now row 38: fix the index number and the updated "want to go" count
$ws1.Range("A38").Value = 37
$ws1.Range("F38").Value = 452

# Update "want to go" counts across sheet 1 rows (scraper refresh)
$ws1.Range("F2").Value = 898
$ws1.Range("F4").Value = 4305
$ws1.Range("F6").Value = 417
$ws1.Range("F7").Value = 3430
$ws1.Range("F8").Value = 963
$ws1.Range("F9").Value = 160
$ws1.Range("F11").Value = 286
$ws1.Range("F12").Value = 2317
$ws1.Range("F13").Value = 1247
$ws1.Range("F14").Value = 26
$ws1.Range("F17").Value = 243
$ws1.Range("F18").Value = 48
$ws1.Range("F19").Value = 9619
$ws1.Range("G19").Value = 75
$ws1.Range("F20").Value = 5906
$ws1.Range("F22").Value = 199
$ws1.Range("F23").Value = 802
$ws1.Range("F24").Value = 111
$ws1.Range("F25").Value = 826
$ws1.Range("F26").Value = 3497
$ws1.Range("F29").Value = 447
$ws1.Range("F30").Value = 97
$ws1.Range("F32").Value = 202
$ws1.Range("F33").Value = 4771
$ws1.Range("F35").Value = 1014
$ws1.Range("F36").Value = 126

# ---- Sheet 2: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F12").Value = 121
$ws2.Range("F15").Value = 3512
$ws2.Range("F23").Value = 10

# ---- Sheet 3: 本地生活 (Local Life) ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 8626
$ws3.Range("F3").Value = 411
$ws3.Range("F4").Value = 1504

# ---- Sheet 4: 全部类型 (All Types) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 8626
$ws4.Range("F3").Value = 898
$ws4.Range("F4").Value = 411
$ws4.Range("F5").Value = 1504
$ws4.Range("F7").Value = 4305
$ws4.Range("F9").Value = 417
$ws4.Range("F10").Value = 3430
$ws4.Range("F11").Value = 963
$ws4.Range("F12").Value = 160
$ws4.Range("F14").Value = 286
$ws4.Range("F15").Value = 2317
$ws4.Range("F19").Value = 1247
$ws4.Range("F21").Value = 26
$ws4.Range("F22").Value = 121
$ws4.Range("F24").Value = 243
$ws4.Range("F25").Value = 48
$ws4.Range("F26").Value = 9620
$ws4.Range("G26").Value = 75
$ws4.Range("F27").Value = 3512
$ws4.Range("F30").Value = 199
$ws4.Range("F31").Value = 802
$ws4.Range("F32").Value = 111
$ws4.Range("F33").Value = 826
$ws4.Range("F34").Value = 3497
$ws4.Range("F37").Value = 447
$ws4.Range("F38").Value = 97
$ws4.Range("F41").Value = 202
$ws4.Range("F42").Value = 4771
$ws4.Range("F43").Value = 1015
$ws4.Range("F44").Value = 126
$ws4.Range("F45").Value = 452
$ws4.Range("F47").Value = 10
